$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The report currently lists data rows 7..75 (item #1..#69) followed by a
# totals row (76) and a footer row (77). A new stock item
# ("حنه جلوري سوده 1 كيس") needs to be inserted as the new row 65 (between
# the existing "حزام فقرات" row and "سرنجات 3 سم" row), which pushes every
# following row down by one: old rows 65-75 become 66-76, the totals row
# moves from 76 to 77 (and its grand total grows by the new item's price),
# and the footer moves from 77 to 78.
# ---------------------------------------------------------------------------

# 1) Push the footer row (77) down to 78, preserving its row height.
$ws.Range("A77:Q77").Copy($ws.Range("A78:Q78"))
$ws.Rows.Item(78).RowHeight = $ws.Rows.Item(77).RowHeight

# 2) Push the totals row (76, only P/Q populated) down to 77.
$ws.Range("A77:O77").Clear()
$ws.Range("P76:Q76").Copy($ws.Range("P77:Q77"))
$ws.Range("P76:Q76").Clear()
$ws.Rows.Item(77).RowHeight = 25.5

# 3) Shift the data rows 65-75 down to 66-76 (bottom-up to avoid clobbering).
for ($r = 75; $r -ge 65; $r--) {
    $dstRow = $r + 1
    $ws.Range("A" + $r + ":Q" + $r).Copy($ws.Range("A" + $dstRow + ":Q" + $dstRow))
    $ws.Rows.Item($dstRow).RowHeight = $ws.Rows.Item($r).RowHeight
}

# 4) Fill in the new row 65 for the new item, matching the existing
#    data-row layout (A:B, C:G, H:K, L:M, N:O merged; P and Q standalone).
$ws.Cells.Item(65, 1).Value = 59
$ws.Cells.Item(65, 3).Value = "حنه جلوري سوده 1 كيس"
$ws.Cells.Item(65, 8).Value = "16:0"
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 14).Value = "40.00"
$ws.Cells.Item(65, 16).Value = "40.0000"
$ws.Cells.Item(65, 17).Value = "1:0"
$ws.Rows.Item(65).RowHeight = 24.75

# 5) Renumber the "م" (sequence) column for the shifted data rows 66-76 so it
#    keeps counting 59,60,61,... beneath the newly inserted row 65.
for ($r = 66; $r -le 76; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 7
}

# 6) Bump the grand total (row 77, column P) by the new item's selling price.
$ws.Cells.Item(77, 16).Value = $ws.Cells.Item(77, 16).Value + 40

# 7) Refresh the generated timestamp string in the footer.
$ws.Cells.Item(78, 1).Value = "Saturday, 4 October, 2025 7:38 PM"

Write-Output "done"
